$wb = $excel.ActiveWorkbook
$w = $excel.ActiveWindow
Write-Output $w.Left
$w.Left = 5500
$w.Top = 500
Write-Output $w.Left
Write-Output $w.Top
